$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 195, shifting existing rows 195:256 down to 196:257
$ws.Rows.Item(195).Insert()

# Populate the newly inserted row 195 with the new data record
$ws.Range("A195").Value = 5
$ws.Range("B195").Value = "Macroferia Regional de Talca"
$ws.Range("C195").Value = "Maule"
$ws.Range("D195").Value = 44900
$ws.Range("E195").Value = 7
$ws.Range("F195").Value = 100112024
$ws.Range("G195").Value = "Choclo"
$ws.Range("H195").Value = "Choclero"
$ws.Range("I195").Value = "Primera"
$ws.Range("J195").Value = 10000
$ws.Range("K195").Value = 350
$ws.Range("L195").Value = 350
$ws.Range("M195").Value = 350
$ws.Range("N195").Value = "`$/unidad"
$ws.Range("O195").Value = "Región de O'Higgins"
$ws.Range("P195").Value = 350
$ws.Range("Q195").Value = 1
$ws.Range("R195").Value = "Hortaliza"

# Match the date-format style used by the other rows' Fecha column (D)
$ws.Range("D195").NumberFormat = $ws.Range("D194").NumberFormat
